$d = $word.ActiveDocument

# The "League and cup will need to be separated" paragraph currently holds
# its text in two runs split around the _GoBack bookmark:
#   "League and cup will need to be separated" + bookmark + " into their respective database"
#
# Target shape:
#   Para A: "League and cup will need to be separated into their respective database"
#   Para B: "" (new blank paragraph)
#   Para C: "Selecting by club code is out as there are multiple clubs with the
#            same code. Best to search by club name" + bookmark (now at the end)

$noteText = "Selecting by club code is out as there are multiple clubs with the same code. Best to search by club name"

# 1) Drop the old trailing run's text - it gets folded into the first run
#    instead. At this point in the document it is the only occurrence of
#    this phrase, so a plain Find/Replace is unambiguous.
$d.Content.Find.Execute(" into their respective database", $true, $false, $false, $false, $false, $true, 1, $false, "", 2)

# 2) Grow the first run's text to the merged sentence, then immediately
#    continue with a paragraph break, a throwaway placeholder token, another
#    paragraph break and the new note text. The placeholder keeps the middle
#    paragraph non-empty for the moment - asking Find/Replace to manufacture
#    a *genuinely* empty paragraph directly leaves a stray empty <w:r/> behind,
#    whereas deleting a placeholder character out of an existing run cleans
#    up after itself properly.
$d.Content.Find.Execute("League and cup will need to be separated", $true, $false, $false, $false, $false, $true, 1, $false, "League and cup will need to be separated into their respective database^pZZZPLACEHOLDERZZZ^p$noteText", 2)

# 3) Strip the placeholder token, leaving a clean, truly empty paragraph.
$d.Content.Find.Execute("ZZZPLACEHOLDERZZZ", $true, $false, $false, $false, $false, $true, 1, $false, "", 2)

# 4) Find the freshly created note paragraph.
$notePara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs($i).Range.Text.StartsWith($noteText)) {
        $notePara = $d.Paragraphs($i)
    }
}

# 5) Relocate the _GoBack bookmark from the middle of the first paragraph to
#    the end of the note paragraph's text. A zero-length range sitting
#    exactly on a paragraph's end-of-text boundary confuses Bookmarks.Add in
#    this runtime, so a throwaway trailing character is inserted first, the
#    bookmark is anchored right before it, and the character is then deleted.
$bm = $d.Bookmarks("_GoBack")
$bm.Delete()

$endPos = $notePara.Range.Start + $noteText.Length
$endRange = $d.Range($endPos, $endPos)
$endRange.InsertAfter("~")

$bmRange = $d.Range($endPos, $endPos)
$d.Bookmarks.Add("_GoBack", $bmRange)

$placeholderChar = $d.Range($endPos, $endPos + 1)
$placeholderChar.Text = ""
